# Update the cryptos list with latest scraped price/volume data.
# Generated on Sat Jan 20 23:47:08 UTC 2024 with GitHub Actions
#
# Helper: writes a value into a cell while forcing a text number format
# first for values that look numeric, so Excel does not silently convert
# the scraped price strings (e.g. "316.23") into floating point numbers
# (which would introduce binary rounding noise and change the cell type).
function Set-TextValue {
    param(
        $Range,
        [string]$Text
    )
    $Range.NumberFormat = "@"
    $Range.Value = $Text
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "41.685.24"
$ws.Range("E2").Value = "  +0.38%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "2.471.22"
$ws.Range("E3").Value = "  -0.81%  "

# Row 4 - TetherUSD
Set-TextValue $ws.Range("D4") "0.999"
$ws.Range("E4").Value = "  +0.09%  "

# Row 5 - BNB
Set-TextValue $ws.Range("D5") "316.23"
$ws.Range("E5").Value = "  +0.66%  "

# Row 6 - Solana
Set-TextValue $ws.Range("D6") "92.74"
$ws.Range("E6").Value = "  -0.68%  "

# Row 7 - XRP
Set-TextValue $ws.Range("D7") "0.554"
$ws.Range("E7").Value = "  +1.78%  "

# Row 8 - USDC
$ws.Range("E8").Value = "  +0.08%  "

# Row 9 - Cardano
Set-TextValue $ws.Range("D9") "0.516"
$ws.Range("E9").Value = "  +2.99%  "

# Row 10 - Dogecoin
Set-TextValue $ws.Range("D10") "0.0882"
$ws.Range("E10").Value = "  +12.14%  "

# Row 11 - Avalanche
Set-TextValue $ws.Range("D11") "32.86"
$ws.Range("E11").Value = "  +0.49%  "

# Row 12 - TRON
$ws.Range("E12").Value = "  -0.48%  "

# Row 13 - WrappedliquidstakedEther2.0
$ws.Range("D13").Value = "2.852.50"
$ws.Range("E13").Value = "  -0.82%  "

# Row 14 - Polkadot
Set-TextValue $ws.Range("D14") "6.91"
$ws.Range("E14").Value = "  +0.63%  "

# Row 15 - Chainlink
Set-TextValue $ws.Range("D15") "15.77"
$ws.Range("E15").Value = "  -2.75%  "

# Row 16 - WrappedEther
$ws.Range("D16").Value = "2.472.22"
$ws.Range("E16").Value = "  +0.65%  "

# Row 17 - Polygon
Set-TextValue $ws.Range("D17") "0.785"
$ws.Range("E17").Value = "  +3.48%  "

# Row 18 - WrappedBTC
$ws.Range("D18").Value = "41.652.16"
$ws.Range("E18").Value = "  +0.23%  "

# Row 19 - ShibaInu
$ws.Range("D19").Value = "0.0₃0972"
$ws.Range("E19").Value = "  +4.48%  "

# Row 20 - Uniswap
$ws.Range("E20").Value = "  +1.97%  "

# Row 21 - Litecoin
Set-TextValue $ws.Range("D21") "71.34"
$ws.Range("E21").Value = "  +0.33%  "

# Row 22 - InternetComputer(DFINITY)
Set-TextValue $ws.Range("D22") "11.46"
$ws.Range("E22").Value = "  +1.86%  "

# Row 23 - BitcoinCash
Set-TextValue $ws.Range("D23") "239.03"
$ws.Range("E23").Value = "  +1.04%  "

# Row 24 - PancakeSwap
$ws.Range("E24").Value = "  +0.38%  "

# Row 25 - ImmutableX
$ws.Range("E25").Value = "  -0.96%  "

# Row 26 - Dai
$ws.Range("E26").Value = "  -0.02%  "

# Row 27 - EthereumClassic
$ws.Range("E27").Value = "  -2.16%  "

# Row 28 - Toncoin
$ws.Range("E28").Value = "  +2.25%  "

# Row 29 - Cosmos
$ws.Range("E29").Value = "  +1.33%  "

# Row 30 - InjectiveProtocol
Set-TextValue $ws.Range("D30") "35.47"
$ws.Range("E30").Value = "  -2.00%  "

# Row 31 - Monero
Set-TextValue $ws.Range("D31") "156.09"
$ws.Range("E31").Value = "  -1.04%  "

# Row 32 - Filecoin
Set-TextValue $ws.Range("D32") "5.53"
$ws.Range("E32").Value = "  +1.04%  "

# Row 33 & 34 - swap Hedera/WEMIXToken ordering
$ws.Range("B33").Value = "WEMIXToken"
$ws.Range("C33").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
Set-TextValue $ws.Range("D33") "2.58"
$ws.Range("E33").Value = "  +0.09%  "

$ws.Range("B34").Value = "Hedera"
$ws.Range("C34").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextValue $ws.Range("D34") "0.0766"
$ws.Range("E34").Value = "  +1.02%  "

# Row 35 - ApeXProtocol
$ws.Range("E35").Value = "  +2.12%  "

# Row 36 - Celestia
Set-TextValue $ws.Range("D36") "17.57"
$ws.Range("E36").Value = "  -1.31%  "

# Row 37 - LidoDAOToken
$ws.Range("E37").Value = "  -2.37%  "

# Row 38 - Stellar
Set-TextValue $ws.Range("D38") "0.115"
$ws.Range("E38").Value = "  +1.29%  "

# Row 39 - ARBITRUM
$ws.Range("E39").Value = "  -2.34%  "

# Row 40 - Kaspa
$ws.Range("E40").Value = "  -2.78%  "

# Row 41 - RenderToken
Set-TextValue $ws.Range("D41") "4.00"
$ws.Range("E41").Value = "  -3.33%  "

# Row 42 - FirstDigitalUSD
$ws.Range("E42").Value = "  +0.02%  "

# Row 43 - Maker
$ws.Range("D43").Value = "1.970.01"
$ws.Range("E43").Value = "  +0.51%  "

# Row 44 - VeChain
$ws.Range("E44").Value = "  +0.13%  "

# Row 45 - EnergySwap
Set-TextValue $ws.Range("D45") "18.94"
$ws.Range("E45").Value = "  -5.14%  "

# Row 46 - NEARProtocol
Set-TextValue $ws.Range("D46") "2.96"
$ws.Range("E46").Value = "  -1.51%  "

# Row 47 - FraxShare
Set-TextValue $ws.Range("D47") "9.08"
$ws.Range("E47").Value = "  +1.92%  "

# Row 48 - RocketPoolETH
$ws.Range("D48").Value = "2.706.27"
$ws.Range("E48").Value = "  -0.83%  "

# Row 49 - Aave
Set-TextValue $ws.Range("D49") "97.45"
$ws.Range("E49").Value = "  +0.58%  "

# Row 50 - ordi
Set-TextValue $ws.Range("D50") "66.97"
$ws.Range("E50").Value = "  -1.36%  "

# Row 51 - MultiversX
Set-TextValue $ws.Range("D51") "52.62"
$ws.Range("E51").Value = "  +3.80%  "
